$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D32").Value = 45906
$ws.Range("E32").Value = 95
$ws.Range("F32").Value = 423
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 1012
$ws.Range("J32").Value = "N/A"

$ws.Range("I46").Select()
